# Auto-generated Excel COM-interop edit script
# Refreshes currentAveragePrice / LevePrice / LeveProfit columns (H:N)
# for a set of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW worksheets,
# matching the upstream scheduled data-refresh commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2110
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 2244.1667
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 6732.500100000001
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -7068.500100000001
# Row 18
$ws.Range("H18").Value = 1041.6666
$ws.Range("I18").Value = 1041.6666
$ws.Range("K18").Value = 1041.6666
$ws.Range("M18").Value = -757.6666
# Row 48
$ws.Range("H48").Value = 10000
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = $null
# Row 56
$ws.Range("H56").Value = 10000
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = $null
# Row 76
$ws.Range("H76").Value = 6000
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = $null
# Row 79
$ws.Range("H79").Value = 6000
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = $null
# Row 125
$ws.Range("H125").Value = 792.26086
$ws.Range("J125").Value = 792.26086
$ws.Range("L125").Value = 7130.34774
$ws.Range("N125").Value = -12050.34774
# Row 127
$ws.Range("H127").Value = 1698.4
$ws.Range("I127").Value = 1623
$ws.Range("K127").Value = 4869
$ws.Range("M127").Value = 91
# Row 132
$ws.Range("H132").Value = 66669344
$ws.Range("I132").Value = 66669344
$ws.Range("K132").Value = 200008032
$ws.Range("M132").Value = -200005502

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1613.6666
$ws.Range("I2").Value = 1336.4
$ws.Range("K2").Value = 1336.4
$ws.Range("M2").Value = -1223.4
# Row 3
$ws.Range("H3").Value = 2750
$ws.Range("I3").Value = 2750
$ws.Range("K3").Value = 2750
$ws.Range("M3").Value = -2635
# Row 4
$ws.Range("H4").Value = 378.33334
$ws.Range("I4").Value = 378.33334
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 378.33334
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -262.33334
$ws.Range("N4").Value = $null
# Row 5
$ws.Range("H5").Value = 622.5
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = $null
# Row 32
$ws.Range("H32").Value = 2748.647
$ws.Range("I32").Value = 2748.647
$ws.Range("K32").Value = 2748.647
$ws.Range("M32").Value = -2461.647
# Row 110
$ws.Range("H110").Value = 15246.706
$ws.Range("I110").Value = 13746.083
$ws.Range("J110").Value = 18848.2
$ws.Range("K110").Value = 13746.083
$ws.Range("L110").Value = 18848.2
$ws.Range("M110").Value = -11701.083
$ws.Range("N110").Value = -22938.2
# Row 116
$ws.Range("H116").Value = 1613.6666
$ws.Range("I116").Value = 1336.4
$ws.Range("K116").Value = 1336.4
$ws.Range("M116").Value = 957.5999999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1613.6666
$ws.Range("I3").Value = 1336.4
$ws.Range("K3").Value = 1336.4
$ws.Range("M3").Value = -1222.4
# Row 4
$ws.Range("H4").Value = 622.5
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null
# Row 11
$ws.Range("H11").Value = 2564.75
$ws.Range("I11").Value = 904
$ws.Range("J11").Value = 3118.3333
$ws.Range("K11").Value = 904
$ws.Range("L11").Value = 3118.3333
$ws.Range("M11").Value = -764
$ws.Range("N11").Value = -3398.3333
# Row 20
$ws.Range("H20").Value = 3649.8
$ws.Range("I20").Value = 3966.3333
$ws.Range("J20").Value = 3175
$ws.Range("K20").Value = 3966.3333
$ws.Range("L20").Value = 3175
$ws.Range("M20").Value = -3719.3333
$ws.Range("N20").Value = -3669

$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 4101.1816
$ws.Range("I19").Value = 510.7
$ws.Range("K19").Value = 510.7
$ws.Range("M19").Value = -340.7
# Row 24
$ws.Range("H24").Value = 4101.1816
$ws.Range("I24").Value = 510.7
$ws.Range("K24").Value = 510.7
$ws.Range("M24").Value = -340.7

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 77.75
$ws.Range("I12").Value = 153.66667
$ws.Range("J12").Value = 52.444443
$ws.Range("K12").Value = 461.00001
$ws.Range("L12").Value = 157.333329
$ws.Range("M12").Value = -288.00001
$ws.Range("N12").Value = -503.333329
# Row 122
$ws.Range("H122").Value = 8292.691999999999
$ws.Range("J122").Value = 692.8570999999999
$ws.Range("L122").Value = 6235.7139
$ws.Range("N122").Value = -11135.7139
# Row 132
$ws.Range("H132").Value = 3990
$ws.Range("I132").Value = 3990
$ws.Range("K132").Value = 35910
$ws.Range("M132").Value = -33380

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 6.4
$ws.Range("I2").Value = 2.4285715
$ws.Range("J2").Value = 15.666667
$ws.Range("K2").Value = 2.4285715
$ws.Range("L2").Value = 15.666667
$ws.Range("M2").Value = 110.5714285
$ws.Range("N2").Value = -241.666667
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = $null
# Row 31
$ws.Range("H31").Value = 1350
$ws.Range("I31").Value = 1350
$ws.Range("K31").Value = 1350
$ws.Range("M31").Value = -1058
# Row 37
$ws.Range("H37").Value = 1350
$ws.Range("I37").Value = 1350
$ws.Range("K37").Value = 1350
$ws.Range("M37").Value = -1073
# Row 70
$ws.Range("H70").Value = 12998.5
$ws.Range("I70").Value = 12998.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 12998.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -12728.5
$ws.Range("N70").Value = $null
# Row 73
$ws.Range("H73").Value = 12998.5
$ws.Range("I73").Value = 12998.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 12998.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -12062.5
$ws.Range("N73").Value = $null
# Row 80
$ws.Range("H80").Value = 5353.1177
$ws.Range("I80").Value = 5967.778
$ws.Range("J80").Value = 4661.625
$ws.Range("K80").Value = 5967.778
$ws.Range("L80").Value = 4661.625
$ws.Range("M80").Value = -4969.778
$ws.Range("N80").Value = -6657.625
# Row 83
$ws.Range("H83").Value = 5353.1177
$ws.Range("I83").Value = 5967.778
$ws.Range("J83").Value = 4661.625
$ws.Range("K83").Value = 29838.89
$ws.Range("L83").Value = 23308.125
$ws.Range("M83").Value = -24846.89
$ws.Range("N83").Value = -33292.125
# Row 122
$ws.Range("H122").Value = 2933.818
$ws.Range("I122").Value = 2933.818
$ws.Range("K122").Value = 8801.454000000002
$ws.Range("M122").Value = -6351.454000000002

$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 5263.375
$ws.Range("I9").Value = 471.8
$ws.Range("J9").Value = 13249.333
$ws.Range("K9").Value = 471.8
$ws.Range("L9").Value = 13249.333
$ws.Range("M9").Value = -247.8
$ws.Range("N9").Value = -13697.333
# Row 82
$ws.Range("H82").Value = 2352.3333
$ws.Range("J82").Value = 3035.5715
$ws.Range("L82").Value = 3035.5715
$ws.Range("N82").Value = -3757.5715
# Row 85
$ws.Range("H85").Value = 2352.3333
$ws.Range("J85").Value = 3035.5715
$ws.Range("L85").Value = 3035.5715
$ws.Range("N85").Value = -5531.5715
# Row 93
$ws.Range("H93").Value = 247.5
$ws.Range("I93").Value = 247.5
$ws.Range("K93").Value = 247.5
$ws.Range("M93").Value = 1000.5
# Row 122
$ws.Range("H122").Value = 4699.5
$ws.Range("I122").Value = 4050
$ws.Range("K122").Value = 12150
$ws.Range("M122").Value = -9700

